$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.326.30"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.878.26"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'0.7226"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.08000"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").Value = "'0.3156"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").Value = "'25.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "'0.08229"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "1.872.80"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "'94.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").Value = "'5.229"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "'0.7126"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'6.420"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008495"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("D18").Value = "29.343.49"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "'243.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").Value = "'13.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'7.774"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'0.1603"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'162.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.040"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "'18.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").Value = "'4.408"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'4.306"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'1.191"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.05%  "
$ws.Range("D32").Value = "'0.05370"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'1.935"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "'0.7603"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").Value = "'1.178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'2.706"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "1.280.45"
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").Value = "'0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").Value = "'2.750"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'6.447"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "'113.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("D42").Value = "'0.9132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("D43").Value = "'74.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("E44").Value = "  +8.36%  "
$ws.Range("D46").Value = "2.030.40"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "'0.5228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "'1.796"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "'9.500"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "'0.4343"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "'7.103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
